# Apply "first set of test cases" edit to locators.xlsx
# Adds two new locator rows (modalContent, modalContentCloseBtn) to the
# digital_coupons_page sheet, right after the existing "loadedText" row.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("digital_coupons_page")

# Insert two blank rows at row 18 (pushes the existing rows 18-28 down to 20-30)
$ws.Range("A18:A19").EntireRow.Insert()

# Row 18: modalContent locator
$ws.Range("B18").Value = "//div[contains(@class,'modal-content')]"
$ws.Range("A18").Value = "modalContent"
$ws.Range("C18").Value = "xpath"

# Row 19: modalContentCloseBtn locator
$ws.Range("B19").Value = "//div[contains(@class,'modal-content')]//button[contains(text(),'Close')]"
$ws.Range("A19").Value = "modalContentCloseBtn"
$ws.Range("C19").Value = "xpath"

# Update the selected/active cell to reflect the new end of the table
$ws.Range("F30").Select()
